# Updated symbol list on Wed Dec 14 07:45:43 UTC 2022 with GitHub Actions
# Refresh the cryptocurrency price/label snapshot in column D (Price) and
# the occasional column E (Volume/label) cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    param(
        [string]$CellRef,
        [string]$Text
    )
    # Column D values are plain-text numeric strings (e.g. "0.005690") whose
    # exact digits/trailing zeros must be preserved. Force the cell to Text
    # format *before* writing so Excel doesn't silently coerce the literal
    # into a floating point number and lose formatting, then drop the
    # temporary "@" number format again so the cell's style matches its
    # original (unstyled) state - only the stored value/type should change.
    $cell = $ws.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.ClearFormats()
}

Set-PriceText "D2"  "274.86"
Set-PriceText "D3"  "23.01"
Set-PriceText "D4"  "6.307"
Set-PriceText "D5"  "0.06268"
Set-PriceText "D7"  "6.680"
Set-PriceText "D8"  "1.362"
Set-PriceText "D9"  "0.8308"
Set-PriceText "D10" "0.01375"
Set-PriceText "D11" "0.1632"
Set-PriceText "D12" "0.08366"
Set-PriceText "D13" "0.03461"
Set-PriceText "D15" "0.09305"
Set-PriceText "D16" "3.891"
Set-PriceText "D17" "0.001658"
Set-PriceText "D18" "0.04772"
Set-PriceText "D19" "0.006332"

Set-PriceText "D20" "0.005690"
$ws.Range("E20").Value = "19HotbitTokenHTBWorstin24h"

Set-PriceText "D21" "0.001094"
Set-PriceText "D22" "0.0001499"
Set-PriceText "D23" "3.717"
Set-PriceText "D24" "2.322"
Set-PriceText "D26" "0.1241"
Set-PriceText "D41" "0.007033"

Set-PriceText "D43" "0.003501"
$ws.Range("E43").Value = "42CEJICEJI"

Set-PriceText "D44" "0.01222"
Set-PriceText "D45" "0.00006253"
Set-PriceText "D47" "0.8996"
Set-PriceText "D48" "0.03960"
Set-PriceText "D49" "0.00002199"
Set-PriceText "D50" "0.01239"
